# repull data, push all data, mean calculation
# Update the dSF column (column F) with freshly re-pulled values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -9
    3  = -8
    4  = -3
    5  = 0
    6  = 2
    7  = -2
    8  = 9
    9  = -1
    10 = -3
    11 = -1
    12 = -4
    13 = 1
    14 = 1
    15 = -1
    16 = 11
    17 = 0
    18 = 3
    20 = 2
    21 = -1
    22 = -2
    23 = -3
    25 = -4
    26 = -1
    27 = 2
    28 = -3
    29 = 1
    30 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
